$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two brand-new data rows right after the header row (new rows 2 & 3). ---
# This shifts the existing data rows 2-21 down to rows 4-23, matching the diff.
$ws.Rows("2:3").Insert()
# The inserted rows pick up formatting from the row above (header); clear it so the
# new data cells are unstyled, matching the rest of the data rows.
$ws.Rows("2:3").ClearFormats()

$newTopRows = @(
    @(-0.0736092627048492, -0.0381790772080421, 0.0797179117798805),
    @(-1.18019163608551,   -4.37715482711792,   0.3266601860523224)
)

$r = 2
foreach ($row in $newTopRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

# --- Append eight brand-new data rows at the bottom (new rows 24-31). ---
$newBottomRows = @(
    @(-0.5545129179954529,  -0.7066183686256409,  -0.1945605874061584),
    @(-0.0233655963093042,  -0.0335975885391235,  -0.5940664410591125),
    @(0.1398881375789642,    0.0471893399953842,   0.531147301197052),
    @(-0.0551305897533893,   0.0639881342649459,   0.093156948685646),
    @(0.1805106848478317,    0.0415388382971286,   0.1635591685771942),
    @(-0.1348485052585601,   0.1539380401372909,   0.1916589736938476),
    @(-0.0897971913218498,   0.1873829066753387,   -0.0282525178045034),
    @(-0.0161879286170005,   0.0589484944939613,   0.0539088584482669)
)

$r = 24
foreach ($row in $newBottomRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}

Write-Output "Applied gyroscope row insert/append edit."
